$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2, header in row 1)
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

# Column C ("Förändrad") holds a date serial value of 45181 for every data row;
# bump it by one day (45182) for all rows from 2 through the last row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
